$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Duplicate row 2's formatting down into the new row 3
$ws.Range("A2:F2").Copy() | Out-Null
$ws.Range("A3:F3").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = $false

# New test case row
$ws.Range("A3").Value = "FUL_Transmittals_ActionRequired_New_ChangeNote"
$ws.Range("B3").Value = "Creates a new Transmittal of  Change Note and validate the count in Action Require"
$ws.Range("C3").Value = "N"
$ws.Range("D3").Value = "Y"
$ws.Range("F3").Value = "Sprint1"

# Extend the list validations down to the new row
$ws.Range("C2:D2").Validation.Delete() | Out-Null
$ws.Range("C2:D3").Validation.Add(3, 1, 1, """Y,N""") | Out-Null

$ws.Range("F2").Validation.Delete() | Out-Null
$ws.Range("F2:F3").Validation.Add(3, 1, 1, """Sprint1,Sprint2,Sprint3,Sprint4,Sprint5,Sprint6,Sprint7,Sprint8,Sprint9,Sprint10""") | Out-Null

$ws.Range("E2:E3").Select()
